$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.Value = "'71.716.68"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +4.60%  "

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.Value = "'4.043.22"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +4.68%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.Value = "'529.82"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.78%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.Value = "'154.47"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +9.39%  "

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.695"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +14.19%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.04%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.760"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.90%  "

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.175"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +5.48%  "

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.0000332"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.83%  "

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.Value = "'49.30"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +18.46%  "

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.Value = "'10.99"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +6.78%  "

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.Value = "'4.694.56"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +4.90%  "

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.Value = "'4.044.67"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +4.84%  "

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.Value = "'14.40"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.89%  "

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.Value = "'21.05"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.11%  "

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.Value = "'1.22"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.90%  "

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.Value = "'0.133"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.19%  "

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.Value = "'71.818.34"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +4.73%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.Value = "'436.59"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +4.69%  "

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.Value = "'3.70"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +6.87%  "

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.Value = "'99.67"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +14.84%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Cells.Item(24, 4)
$c.Value = "'14.77"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +5.69%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "PancakeSwap"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Cells.Item(25, 4)
$c.Value = "'4.25"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +6.57%  "

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.Value = "'11.41"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.35%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.Value = "'10.92"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +3.66%  "

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.Value = "'37.22"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.18%  "

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.Value = "'5.84"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +3.01%  "

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.Value = "'3.53"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +25.54%  "

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.Value = "'13.66"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.70%  "

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.132"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +5.89%  "

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.Value = "'673.07"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.52%  "

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.Value = "'6.71"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.58%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.Value = "'66.87"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.31%  "

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.Value = "'42.72"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +7.69%  "

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.435"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.68%  "

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.0₃0859"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.10%  "

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.157"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +6.26%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.Value = "'3.42"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.99%  "

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.05%  "

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.0499"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.25%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.06%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +3.19%  "

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.152"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +8.77%  "

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.Value = "'2.74"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.80%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Cells.Item(47, 4)
$c.Value = "'3.38"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.52%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "THORChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'9.52"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +10.85%  "

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.Value = "'3.06"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.79%  "

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.000274"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.36%  "

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.Value = "'3.36"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.87%  "
